# Convert a "RRGGBB" hex string into the integer encoding used by the
# PowerPoint COM RGB()-style color properties (r + g*256 + b*65536).
function HexToRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# --- 1. Re-colour the deck's theme (the one backing every slide master) ---
# from the "Integral" / "Red Violet" palette back to the stock
# "Office" palette. The 12 slots follow clrScheme order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$s1 = $p.Slides.Item(1)
$tcs = $s1.ThemeColorScheme
for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Colors($i).RGB = HexToRgb $officeColors[$i - 1]
}

# --- 2. Re-apply the built-in table style on the B1 financial-documents
# table (slide 5) ---
$s5 = $p.Slides.Item(5)
for ($i = 1; $i -le $s5.Shapes.Count; $i++) {
    $shp = $s5.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{C302143A-D234-4C31-A900-6FBE34561C39}")
    }
}
